# Rename the i18n locale-header columns from the bare locale codes
# ("en"/"ua") to the fully-qualified "value.en" / "value.ua" keys, and
# make the "i18n" sheet the active/selected sheet+cell (it was
# "metadata" before), matching the edit that ships with the
# table-generator/controllers work.

$wb = $excel.ActiveWorkbook

$i18n = $wb.Worksheets.Item("i18n")
$i18n.Range("B1").Value = "value.en"
$i18n.Range("C1").Value = "value.ua"

# Activate the i18n sheet and move the selection to C11, mirroring the
# author's last-saved cursor position.
$i18n.Activate()
$i18n.Range("C11").Select()
